# Update mapping suite metadata + fix entryRelationship/extension/reasonReference
# path separators (dot -> colon) per NRISS review.

$wb = $excel.ActiveWorkbook

# ----- Metadata sheet -----
$meta = $wb.Worksheets.Item("Metadata")

# The "Name" row value is cleared; its text moves up to become the "Title" row value.
$meta.Range("B4").ClearContents()
$meta.Range("B5").Value = "Mapping Métier/CDA/FHIR : `"Prescription de dispositif médical`""

# Refresh the Date row value.
$meta.Range("B8").Value = "2026-01-07T15:20:53+00:00"

# ----- Mapping Table 0 (FRLM -> FRCDA) -----
$map0 = $wb.Worksheets.Item("Mapping Table 0")
$map0.Range("D12").Value = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecALD"
$map0.Range("D13").Value = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecAccidentTravail"
$map0.Range("D14").Value = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecPrevention"
$map0.Range("D15").Value = "FRCDADispositifMedical.entryRelationship:frNonRemboursable"

# ----- Mapping Table 1 (FRCDA -> FRDeviceRequestDocument) -----
$map1 = $wb.Worksheets.Item("Mapping Table 1")
$map1.Range("D10").Value = "FRDeviceRequestDocument.requester.extension:prescripteur"
$map1.Range("A12").Value = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecALD"
$map1.Range("D12").Value = "FRDeviceRequestDocument.reasonReference:EnRapportAvecALD"
$map1.Range("A13").Value = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecAccidentTravail"
$map1.Range("D13").Value = "FRDeviceRequestDocument.reasonReference:EnRapportAvecAccidentTravail"
$map1.Range("A14").Value = "FRCDADispositifMedical.entryRelationship:frEnRapportAvecPrevention"
$map1.Range("D14").Value = "FRDeviceRequestDocument.reasonReference:EnRapportAvecLaPrevention"
$map1.Range("A15").Value = "FRCDADispositifMedical.entryRelationship:frNonRemboursable"
$map1.Range("D15").Value = "FRDeviceRequestDocument.extension:notCovered"
